$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for storage pool metabolites: L-Valine (row 10) and L-Aspartate (row 17)
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(17).Insert()

# Apply the same label style (bold, thin border, centered) used by column A to the new rows
$newLabelRows = @(10, 17)
foreach ($r in $newLabelRows) {
    $cell = $ws.Range("A" + $r)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Write final recomputed storage-pool metabolite table (values refreshed after fixing storage pool metabolites)
$ws.Range("A2").Value = "Sucrose_c_Day_sp_exchange"
$ws.Range("B2").Value = 0.02497
$ws.Range("C2").Value = -0.2960924999994341
$ws.Range("D2").Value = 0.2960924999994431
$ws.Range("E2").Value = $false
$ws.Range("A3").Value = "L-Isoleucine_Day_sp_exchange"
$ws.Range("B3").Value = 0.0033
$ws.Range("C3").Value = -0.0033
$ws.Range("D3").Value = 0.0033
$ws.Range("E3").Value = $false
$ws.Range("A4").Value = "L-Leucine_Day_sp_exchange"
$ws.Range("B4").Value = 0.0077
$ws.Range("C4").Value = -0.0077
$ws.Range("D4").Value = 0.0077
$ws.Range("E4").Value = $false
$ws.Range("A5").Value = "L-Lysine_Day_sp_exchange"
$ws.Range("B5").Value = 0.00005500000000000028
$ws.Range("C5").Value = -0.000055
$ws.Range("D5").Value = 0.000055
$ws.Range("E5").Value = $false
$ws.Range("A6").Value = "L-Methionine_Day_sp_exchange"
$ws.Range("B6").Value = 0.00011
$ws.Range("C6").Value = -0.0001100000000000026
$ws.Range("D6").Value = 0.0001100000000000028
$ws.Range("E6").Value = $false
$ws.Range("A7").Value = "L-Phenylalanine_Day_sp_exchange"
$ws.Range("B7").Value = 0.004400000000000007
$ws.Range("C7").Value = -0.02794000000000093
$ws.Range("D7").Value = 0.02794000000000001
$ws.Range("E7").Value = $false
$ws.Range("A8").Value = "L-Threonine_Day_sp_exchange"
$ws.Range("B8").Value = 0.004729999999999999
$ws.Range("C8").Value = -0.01154999999999798
$ws.Range("D8").Value = 0.01155000000000001
$ws.Range("E8").Value = $false
$ws.Range("A9").Value = "L-Tryptophan_Day_sp_exchange"
$ws.Range("B9").Value = 0.002200000000000001
$ws.Range("C9").Value = -0.002200000000000001
$ws.Range("D9").Value = 0.002200000000000001
$ws.Range("E9").Value = $false
$ws.Range("A10").Value = "L-Valine_Day_sp_exchange"
$ws.Range("B10").Value = 0.005500000000000002
$ws.Range("C10").Value = -0.0209
$ws.Range("D10").Value = 0.0209
$ws.Range("E10").Value = $false
$ws.Range("A11").Value = "L-Cysteine_Day_sp_exchange"
$ws.Range("B11").Value = 0.00132
$ws.Range("C11").Value = -0.00131999999981512
$ws.Range("D11").Value = 0.001319999999866685
$ws.Range("E11").Value = $false
$ws.Range("A12").Value = "L-Glutamine_c_Day_sp_exchange"
$ws.Range("B12").Value = 0.02048200000000043
$ws.Range("C12").Value = -0.1926048928566449
$ws.Range("D12").Value = 0.481076749999466
$ws.Range("E12").Value = $false
$ws.Range("A13").Value = "L-Glutamate_c_Day_sp_exchange"
$ws.Range("B13").Value = 1.416330230160195
$ws.Range("C13").Value = 1.203238055557414
$ws.Range("D13").Value = 7.611845999993817
$ws.Range("E13").Value = $true
$ws.Range("A14").Value = "L-Tyrosine_Day_sp_exchange"
$ws.Range("B14").Value = 0.009569999999999997
$ws.Range("C14").Value = -0.009569999999995785
$ws.Range("D14").Value = 0.009569999999983739
$ws.Range("E14").Value = $false
$ws.Range("A15").Value = "L-Asparagine_Day_sp_exchange"
$ws.Range("B15").Value = 0.02145
$ws.Range("C15").Value = -0.1488162499994153
$ws.Range("D15").Value = 0.1488162499996062
$ws.Range("E15").Value = $false
$ws.Range("A16").Value = "L-Serine_c_Day_sp_exchange"
$ws.Range("B16").Value = 0.01155
$ws.Range("C16").Value = -0.2733224999990567
$ws.Range("D16").Value = 0.2733224999988936
$ws.Range("E16").Value = $false
$ws.Range("A17").Value = "L-Aspartate_c_Day_sp_exchange"
$ws.Range("B17").Value = 0.2722904206341817
$ws.Range("C17").Value = -4.424533666663168
$ws.Range("D17").Value = 0.615861888888184
$ws.Range("E17").Value = $false
$ws.Range("A18").Value = "Starch_p_Day_sp_exchange"
$ws.Range("B18").Value = 0.1141292777775617
$ws.Range("C18").Value = -0.2856974999999945
$ws.Range("D18").Value = 0.2609575833331929
$ws.Range("E18").Value = $false
$ws.Range("A19").Value = "(S)-Malate_c_Day_sp_exchange"
$ws.Range("B19").Value = 1.073623309526004
$ws.Range("C19").Value = 0.99116783333576
$ws.Range("D19").Value = 9.399136999984796
$ws.Range("E19").Value = $true
$ws.Range("A20").Value = "Fumarate_Day_sp_exchange"
$ws.Range("B20").Value = 0.0009240000000030056
$ws.Range("C20").Value = -0.03814799999993534
$ws.Range("D20").Value = 0.03814799999988958
$ws.Range("E20").Value = $false
